# 822-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-Late Repayment-Makerepayment1.xlsx
# Commit: "Loan RBI, Variable Instalments"
#
# The "Repayment Schedule" sheet gets a new column inserted before column O
# (pushing the old O/P data right into P/Q), and the data that used to live
# in column N is moved into the freshly inserted column O, leaving N blank.
# The active sheet/tab also switches from "Transactions" back to
# "Repayment Schedule", with a new selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column at O - shifts old O->P and P->Q.
$ws.Columns("O").Insert()

# Move what used to be in column N into the newly-opened column O, then
# blank out N (keeping the cell/style in place).
for ($r = 1; $r -le 15; $r++) {
    $nCell = $ws.Cells.Item($r, 14)
    $oCell = $ws.Cells.Item($r, 15)
    $oCell.Value = $nCell.Value()
    $nCell.ClearContents()
}

# Update the selection on the Repayment Schedule sheet and make it active.
$ws.Range("R8").Select()
$ws.Activate()
